$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Row=2;  B=0.7287194209349384;  C=0.004309184025731883; D=0.1529057820181812; E=0.4998867070740569; G=1.385821094052908},
    @{Row=3;  B=0.1554434735375247;  C=0.3375848360084654;   D=3.082599426703578;  E=0.4998867070740569; G=4.075514443323626},
    @{Row=4;  B=3.182878228561681;   C=1.65323645889881;     D=3.082599426703578;  E=0.4998867070740569; G=8.418600821238126},
    @{Row=5;  B=3.182878228561681;   C=1.65323645889881;     D=3.082599426703578;  E=0.4998867070740569; G=8.418600821238126},
    @{Row=6;  B=1.505614041169197;   C=1.65323645889881;     D=3.082599426703578;  E=6.48142807727062;   G=12.7228780040422},
    @{Row=7;  B=3.182878228561681;   C=1.65323645889881;     D=0.7127328510149897; E=0.4998867070740569; G=6.048734245549538},
    @{Row=8;  B=1.505614041169197;   C=9.226618575922256;    D=1935279062.313128;  E=71517.89157740913;  G=1935350590.936939},
    @{Row=9;  B=3.182878228561681;   C=1.65323645889881;     D=3.082599426703578;  E=6.48142807727062;   G=14.40014219143469},
    @{Row=10; B=0.06328177979961902; C=0.004309184025731883; D=0.1529057820181812; E=6.48142807727062;   G=6.701924823114153},
    @{Row=11; B=3.182878228561681;   C=1.65323645889881;     D=3.082599426703578;  E=0.4998867070740569; G=8.418600821238126}
)

foreach ($r in $data) {
    $row = $r.Row
    $ws.Range("B$row").Value = $r.B
    $ws.Range("C$row").Value = $r.C
    $ws.Range("D$row").Value = $r.D
    $ws.Range("E$row").Value = $r.E
    $ws.Range("G$row").Value = $r.G
}
